$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell C11 from 3 to 29
$ws.Range("C11").Value = 29

# Update the active selection to D5
$ws.Range("D5").Select()
